$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44187
$ws.Range("K3").Value = "Dina"
$ws.Range("M3").Value = 55
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 16000
$ws.Range("P3").Value = 15455
$ws.Range("Q3").Value = '$/caja 15 kilos granel'
$ws.Range("S3").Value = 1030

# Row 4
$ws.Range("D4").Value = 44174
$ws.Range("K4").Value = "Castle Brite"
$ws.Range("M4").Value = 75
$ws.Range("N4").Value = 9000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 9467
$ws.Range("Q4").Value = '$/caja 10 kilos'
$ws.Range("S4").Value = 947
$ws.Range("T4").Value = 10

# Row 5
$ws.Range("D5").Value = 44176
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 17000
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 17400
$ws.Range("Q5").Value = '$/caja 18 kilos'
$ws.Range("S5").Value = 967
$ws.Range("T5").Value = 18

# Row 8
$ws.Range("D8").Value = 44537
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 21000
$ws.Range("O8").Value = 21500
$ws.Range("P8").Value = 21250
$ws.Range("Q8").Value = '$/caja 15 kilos'
$ws.Range("S8").Value = 1417
$ws.Range("T8").Value = 15

# Row 9
$ws.Range("D9").Value = 44181
$ws.Range("K9").Value = "Modesto"
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 21000
$ws.Range("P9").Value = 20500
$ws.Range("Q9").Value = '$/caja 18 kilos'
$ws.Range("R9").Value = "Región de Coquimbo"
$ws.Range("S9").Value = 1139
$ws.Range("T9").Value = 18

# Row 10
$ws.Range("D10").Value = 44552
$ws.Range("K10").Value = "Castle Brite"
$ws.Range("M10").Value = 120
$ws.Range("N10").Value = 15500
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 15750
$ws.Range("Q10").Value = '$/caja 15 kilos'
$ws.Range("R10").Value = "Región de O'Higgins"
$ws.Range("S10").Value = 1050
$ws.Range("T10").Value = 15
